$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell (H1) onto the two new header
# cells so they reuse the same style (bold, bordered, centered) instead of
# minting a brand new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
